$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("20190221_run")

# --- Insert the two new EAG rows (bottom-most first so row numbers stay stable) ---
$ws.Rows.Item(13).Insert()   # room for 3110-EAG-1, just above the old row 13 (3110-EAG-3)
$ws.Range("A13:B13").Style = "Normal"   # inserted row inherits formatting from above; reset it

$ws.Rows.Item(2).Insert()    # room for 2130-EAG-2, right after the header row
$ws.Range("A2:B2").Style = "Normal"

# --- Fill in column A (EAG id) and column B (Opmerking) for every data row ---
$ws.Cells.Item(2,1).Value = '2130-EAG-2'
$ws.Cells.Item(2,2).Value = 'q_cso iets anders dan in Excel. Poel en dijkslek inlaat zijn eerst ValueSeries en daarna handmatig ingevoerde reeks. ValueSeries ontbreekt nu uit reeksen.csv.'
$ws.Cells.Item(3,1).Value = '2140-EAG-3'
$ws.Cells.Item(3,2).Value = 'intrek heeft een minimum in Excel die niet in Python staat. Peil, inlaat en uitlaat waarshcijnlijk om die reden ook fout.'
$ws.Cells.Item(4,1).Value = '2140-EAG-6'
$ws.Cells.Item(4,2).Value = 'intrek en drain niet juist. q_cso ook iets anders maar ws minder groot probleem. Hopelijk alle verschillen functie van intrek fout'
$ws.Cells.Item(5,1).Value = '2250-EAG-2'
$ws.Cells.Item(5,2).Value = 'peil totaal anders, ontbrekende ValueSeries lijkt het bij verschil in drain en berekende uitlaat. Eerst peil beschouwen, dan verder kijken.'
$ws.Cells.Item(6,1).Value = '2500-EAG-6'
$ws.Cells.Item(6,2).Value = 'verschil intrek, daardoor verschil in inlaat, uitlaat en peil? Hoge piek in uitspoeling aan begin, waardoor?'
$ws.Cells.Item(7,1).Value = '2501-EAG-1'
$ws.Cells.Item(7,2).Value = 'Teveel inlaat, mogelijk iets met peil hTargets? Ook uitlaat klopt niet maar ws als gevolg van iets anders.'
$ws.Cells.Item(8,1).Value = '2501-EAG-2'
$ws.Cells.Item(8,2).Value = 'Goed'
$ws.Cells.Item(9,1).Value = '2505-EAG-1'
$ws.Cells.Item(9,2).Value = 'Goed'
$ws.Cells.Item(10,1).Value = '2510-EAG-2'
$ws.Cells.Item(10,2).Value = 'Goed'
$ws.Cells.Item(11,1).Value = '2510-EAG-3'
$ws.Cells.Item(11,2).Value = 'Goed'
$ws.Cells.Item(12,1).Value = '3050-EAG-1'
$ws.Cells.Item(12,2).Value = 'Verschil uitspoeling, Python heeft minimum grens, Excel niet. Verschil q_cso.'
$ws.Cells.Item(13,1).Value = '3050-EAG-2'
$ws.Cells.Item(13,2).Value = 'Verschil uitspoeling, missende ValueSeries lijkt het in verschil in uitlaat. Q_cso anders. Geen inlaat in Excel, wel in Python af en toe.'
$ws.Cells.Item(14,1).Value = '3110-EAG-1'
$ws.Cells.Item(14,2).Value = 'Later kwel en wegzijging reeksen in Excel niet in Python ingevoerd. Beginfase verschil door peil? In excel niet onder ondergrens, in python wel.'
$ws.Cells.Item(15,1).Value = '3110-EAG-3'
$ws.Cells.Item(15,2).Value = 'Goed'
$ws.Cells.Item(16,1).Value = '3200-EAG-2'
$ws.Cells.Item(16,2).Value = 'Goed (minimaal piekje bij begin intrek)'
$ws.Cells.Item(17,1).Value = '3201-EAG-1'
$ws.Cells.Item(17,2).Value = 'Iets fout met peil berekening'
$ws.Cells.Item(18,1).Value = '3201-EAG-2'
$ws.Cells.Item(18,2).Value = 'Iets fout met peil berekening'
$ws.Cells.Item(19,1).Value = '3201-EAG-3'
$ws.Cells.Item(19,2).Value = 'Verschil intrek, verschil uitspoeling, verschil peil.'
$ws.Cells.Item(20,1).Value = '3230-EAG-1'
$ws.Cells.Item(20,2).Value = 'Verschil uitspoeling (systematisch), verschil q_cso, klein verschil peil, veel te veel uitlaat (geen uitlaat in Excel)'
$ws.Cells.Item(21,1).Value = '3230-EAG-2'
$ws.Cells.Item(21,2).Value = 'Goed (verschil millimeters in peil, teveel inlaat en uitlaat, maar verschil mogelijk veroorzaakt door verdamping?)'
$ws.Cells.Item(22,1).Value = '3230-EAG-3'
$ws.Cells.Item(22,2).Value = 'Goed'
$ws.Cells.Item(23,1).Value = '3230-EAG-4'
$ws.Cells.Item(23,2).Value = 'Verschil uitspoeling (systematisch), verschil q_cso, peil ook niet goed.'
$ws.Cells.Item(24,1).Value = '3230-EAG-5'
$ws.Cells.Item(24,2).Value = 'Goed? (minimaal systematisch verschil uitspoeling, verschil uitlaat en inlaat door verdamping?)'
$ws.Cells.Item(25,1).Value = '3260-EAG-1'
$ws.Cells.Item(25,2).Value = 'Goed'
$ws.Cells.Item(26,1).Value = '3301-EAG-1'
$ws.Cells.Item(26,2).Value = 'intrek verschil, piek uitspoeling aan begin (kleine fout)'
$ws.Cells.Item(27,1).Value = '3301-EAG-2'
$ws.Cells.Item(27,2).Value = 'intrek verschil, piek uitspoeling aan begin (kleine fout)'
$ws.Cells.Item(28,1).Value = '3303-EAG-1'
$ws.Cells.Item(28,2).Value = 'intrek verschil, piek uitspoeling aan begin (kleine fout)'
$ws.Cells.Item(29,1).Value = '3360-EAG-1'
$ws.Cells.Item(29,2).Value = 'Goed'

# --- Apply the built-in "Good" cell style (green fill/text) to rows with a positive verdict ---
$goodRows = @(8,9,10,11,15,16,21,22,24,25,29)
foreach ($r in $goodRows) {
    $ws.Cells.Item($r,1).Style = "Good"
    $ws.Cells.Item($r,2).Style = "Good"
}

# --- Column B needs to be wide enough for the longer comments now stored in it ---
$ws.Columns.Item(2).ColumnWidth = 124.9296875

# --- Match the saved selection/active cell from the authored workbook ---
$ws.Range("B27").Select()